$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.312.27'
$ws.Range('E2').Value = '  +0.30%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.690.52'
$ws.Range('E3').Value = '  +1.27%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.008'
$ws.Range('E4').Value = '  +0.21%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '218.74'
$ws.Range('E5').Value = '  +0.32%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5259'
$ws.Range('E6').Value = '  +3.85%  '
$ws.Range('E7').Value = '  +0.10%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2704'
$ws.Range('E8').Value = '  +1.70%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06431'
$ws.Range('E9').Value = '  +1.45%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '22.05'
$ws.Range('E10').Value = '  +2.24%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07495'
$ws.Range('E11').Value = '  +1.76%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.726.80'
$ws.Range('E12').Value = '  +3.37%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.566'
$ws.Range('E13').Value = '  +0.66%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.5854'
$ws.Range('E14').Value = '  +0.94%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.000008527'
$ws.Range('E15').Value = '  -0.01%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.59'
$ws.Range('E16').Value = '  -0.34%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '26.367.82'
$ws.Range('E17').Value = '  +0.78%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '4.951'
$ws.Range('E18').Value = '  +0.35%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.008'
$ws.Range('E19').Value = '  +0.15%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '10.90'
$ws.Range('E20').Value = '  +0.71%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '189.80'
$ws.Range('E21').Value = '  +0.27%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.219'
$ws.Range('E22').Value = '  +0.42%  '
$ws.Range('E23').Value = '  +0.13%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '144.55'
$ws.Range('E24').Value = '  +0.51%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '7.700'
$ws.Range('E25').Value = '  +0.11%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1232'
$ws.Range('E26').Value = '  +5.36%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '15.87'
$ws.Range('E27').Value = '  +1.10%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.06671'
$ws.Range('E28').Value = '  +15.45%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.356'
$ws.Range('E29').Value = '  +6.14%  '
$ws.Range('E30').Value = '  +0.45%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.589'
$ws.Range('E31').Value = '  +2.14%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.580'
$ws.Range('E32').Value = '  +1.45%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.671'
$ws.Range('E33').Value = '  +2.28%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.032'
$ws.Range('E34').Value = '  +2.04%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.6238'
$ws.Range('E35').Value = '  +4.43%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.397'
$ws.Range('E36').Value = '  +1.59%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.705'
$ws.Range('E37').Value = '  +2.41%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.388'
$ws.Range('E38').Value = '  +6.26%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.112.28'
$ws.Range('E39').Value = '  +3.60%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01624'
$ws.Range('E40').Value = '  +0.78%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8887'
$ws.Range('E41').Value = '  +3.32%  '
$ws.Range('E42').Value = '  +0.78%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '100.82'
$ws.Range('E43').Value = '  +1.32%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.838.89'
$ws.Range('E44').Value = '  +1.17%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.00000000115'
$ws.Range('E45').Value = '  +3.69%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '57.02'
$ws.Range('E46').Value = '  +2.35%  '
$ws.Range('B47').Value = 'Frax'
$ws.Range('C47').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.011'
$ws.Range('E47').Value = '  +0.86%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.159'
$ws.Range('E48').Value = '  +1.11%  '
$ws.Range('E49').Value = '  +1.72%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.4305'
$ws.Range('E50').Value = '  +0.09%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '6.069'
$ws.Range('E51').Value = '  +3.54%  '

Write-Host "Applied cryptos update"
